# Adds two new columns, I ("I0") and J ("IF"), to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header cell's formatting (bold, centered, thin border) from H1
# onto the two new header cells so they visually match the rest of the
# header row, then set their text.
$ws.Range("H1").Copy($ws.Range("I1:J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for the new I/J columns, keyed by row number.
$data = @(
    @(2, 2, 2),
    @(3, 7, 7),
    @(4, 8, 8),
    @(5, 7, 7),
    @(6, 5, 5),
    @(7, 8, 8),
    @(8, 6, 6),
    @(9, 9, 9),
    @(10, 9, 9),
    @(11, 6, 7),
    @(12, 10, 10),
    @(13, 5, 6),
    @(14, 8, 8),
    @(15, 8, 8),
    @(16, 7, 7),
    @(17, 5, 6),
    @(18, 5, 5),
    @(19, 8, 8),
    @(20, 8, 8),
    @(21, 9, 9),
    @(22, 8, 8),
    @(23, 6, 6),
    @(24, 9, 9),
    @(25, 6, 6),
    @(26, 2, 3),
    @(27, 4, 4),
    @(28, 8, 8),
    @(29, 6, 6)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($row, 9).Value = $iVal
    $ws.Cells.Item($row, 10).Value = $jVal
}

Write-Output "I0/IF columns added"
